{"js": "// Office.js (Word JavaScript API) script.\n// Body of: async (context) => { ... }\n//\n// Change 1: Insert a new empty paragraph (same paragraph formatting:\n// style \"List Paragraph\" / \"Prrafodelista\", right indent -660 twips,\n// justified) immediately before the paragraph that begins with\n// \"EL PAGO DEBER\u00c1 REALIZARSE DE LUNES A S\u00c1BADO...\".\n//\n// Change 2: Inside the clause \"...QUEDAR\u00c1N OBLIGADOS A CUBRIR LA PENA\n// CONVENCIONAL.\" turn \"OBLIGADOS\" into \"OBLIGAD{{SEXO_11}}S\" (a gender\n// merge-field, matching the template's existing {{SEXO_n}} pattern).\n\n// --- Change 1: insert blank paragraph before \"EL PAGO DEBER\u00c1 REALIZARSE...\" ---\nconst payParas = context.document.body.search(\n  \"EL PAGO DEBER\u00c1 REALIZARSE DE LUNES A S\u00c1BADO\",\n  { matchCase: true }\n);\npayParas.load(\"items\");\nawait context.sync();\n\nif (payParas.items.length > 0) {\n  const payParagraph = payParas.items[0].paragraphs.getFirst();\n  payParagraph.insertParagraph(\"\", \"Before\");\n}\n\n// --- Change 2: split \" OBLIGADOS\" to insert the {{SEXO_11}} merge field ---\nconst obligadosRange = context.document.body.search(\n  \"DOS A CUBRIR\",\n  { matchCase: true }\n);\nobligadosRange.load(\"items\");\nawait context.sync();\n\nif (obligadosRange.items.length > 0) {\n  obligadosRange.items[0].insertText(\"D{{SEXO_11}}S A CUBRIR\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument is the open document ($d below).\n#\n# Change 1: Insert a new empty paragraph (inheriting the same paragraph\n# formatting: style \"List Paragraph\" / \"Prrafodelista\", right indent\n# -660 twips, justified) immediately before the paragraph that begins\n# with \"EL PAGO DEBER\u00c1 REALIZARSE DE LUNES A S\u00c1BADO...\".\n#\n# Change 2: Inside the clause \"...QUEDAR\u00c1N OBLIGADOS A CUBRIR LA PENA\n# CONVENCIONAL.\" turn \"OBLIGADOS\" into \"OBLIGAD{{SEXO_11}}S\" (a gender\n# merge-field, matching the template's existing {{SEXO_n}} pattern).\n\n$d = $word.ActiveDocument\n\n# --- Change 1: insert blank paragraph before \"EL PAGO DEBER\u00c1 REALIZARSE...\" ---\n$rng = $d.Content\n$found = $rng.Find.Execute(\"EL PAGO DEBER\u00c1 REALIZARSE DE LUNES A S\u00c1BADO\")\nif ($found) {\n    $targetParagraph = $rng.Paragraphs(1)\n    $targetParagraph.Range.InsertParagraphBefore()\n}\n\n# --- Change 2: split \"OBLIGADOS\" to insert the {{SEXO_11}} merge field ---\n$rng2 = $d.Content\n$found2 = $rng2.Find.Execute(\"OBLIGADOS A CUBRIR LA PENA CONVENCIONAL\")\nif ($found2) {\n    $rng2.Text = \"OBLIGAD{{SEXO_11}}S A CUBRIR LA PENA CONVENCIONAL\"\n}\n"}
